$d = $word.ActiveDocument

# Mapping of old text -> new text (wildcard off, match case, whole match)
$replacements = @(
    @("2025-08-28 Thursday", "2025-08-29 Friday"),
    @("750×6=", "628×6="),
    @("116×9=", "556×5="),
    @("674×4=", "117×2="),
    @("848×9=", "775×6="),
    @("703×6=", "679×3="),
    @("177×8=", "690×9="),
    @("343×6=", "690×8="),
    @("767×7=", "982×9="),
    @("669×3=", "997×7="),
    @("271×2=", "444×8="),
    @("455×3=", "420×7="),
    @("812×7=", "509×5="),
    @("933×9=", "803×9="),
    @("262×9=", "427×4="),
    @("300×9=", "819×2="),
    @("804×9=", "670×9="),
    @("720×3=", "410×8="),
    @("542×6=", "509×6="),
    @("149×6=", "928×9="),
    @("604×8=", "601×3="),
    @("403×7=", "401×8="),
    @("371×9=", "475×5="),
    @("984×7=", "404×8="),
    @("294×6=", "135×8="),
    @("949×5=", "814×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
